{"js": "// Locate the \"GIS & Geospatial Analysis Consulting\" paragraph (under the\n// \"PARTNER - Siege Analytics\" heading) and insert three new bullet\n// paragraphs directly after it, before the existing \"\u2022 Lead comprehensive\n// research initiatives...\" bullet.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"GIS & Geospatial Analysis Consulting\"\n);\n\nif (!target) {\n  throw new Error('Could not find paragraph \"GIS & Geospatial Analysis Consulting\"');\n}\n\nconst newBullets = [\n  \"\\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n  \"\\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n  \"\\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\",\n];\n\n// Insert in order, always right after the target paragraph, so each\n// subsequent insertion lands directly below the previous new bullet\n// (maintaining the original order from the diff).\nlet anchor = target;\nfor (const text of newBullets) {\n  anchor = anchor.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Locate the \"GIS & Geospatial Analysis Consulting\" paragraph (under the\n# \"PARTNER - Siege Analytics\" heading) and insert three new bullet\n# paragraphs directly after it, before the existing \"\u2022 Lead comprehensive\n# research initiatives...\" bullet.\n$d = $word.ActiveDocument\n\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"GIS & Geospatial Analysis Consulting\") {\n        $target = $i\n        break\n    }\n}\n\nif ($target -eq $null) {\n    Write-Output \"Could not find paragraph 'GIS & Geospatial Analysis Consulting'\"\n} else {\n    $bullets = @(\n        \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n        \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n        \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n    )\n\n    $insertAt = $target\n    foreach ($bulletText in $bullets) {\n        $p = $d.Paragraphs.Item($insertAt)\n        $p.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($insertAt + 1)\n        $newPara.Range.Text = $bulletText\n        $insertAt = $insertAt + 1\n    }\n}\n"}
